# Apply updated market-price / profit figures across the Sheets workbook.
# Generated from the authoritative cell-level diff (old -> new values).

$wb = $excel.ActiveWorkbook

# ----- ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 317.73685
$ws.Range("I28").Value = 317.86667
$ws.Range("K28").Value = 317.86667
$ws.Range("M28").Value = 167.13333
$ws.Range("H33").Value = 185.91667
$ws.Range("I33").Value = 185.91667
$ws.Range("K33").Value = 185.91667
$ws.Range("M33").Value = 43.08332999999999
$ws.Range("H53").Value = 89.375
$ws.Range("I53").Value = 100.71429
$ws.Range("K53").Value = 100.71429
$ws.Range("M53").Value = 536.28571
$ws.Range("H76").Value = 4000
$ws.Range("J76").Value = 4000
$ws.Range("L76").Value = 4000
$ws.Range("N76").Value = -4630
$ws.Range("H79").Value = 4000
$ws.Range("J79").Value = 4000
$ws.Range("L79").Value = 4000
$ws.Range("N79").Value = -6184
$ws.Range("H80").Value = 461.33334
$ws.Range("J80").Value = 750
$ws.Range("L80").Value = 2250
$ws.Range("N80").Value = -4246
$ws.Range("H83").Value = 461.33334
$ws.Range("J83").Value = 750
$ws.Range("L83").Value = 6750
$ws.Range("N83").Value = -16734
$ws.Range("H107").Value = 1216.8
$ws.Range("I107").Value = 977.25
$ws.Range("J107").Value = 1490.5714
$ws.Range("K107").Value = 977.25
$ws.Range("L107").Value = 1490.5714
$ws.Range("M107").Value = 942.75
$ws.Range("N107").Value = -5330.5714
$ws.Range("H132").Value = 4129.2144
$ws.Range("I132").Value = 2790.3
$ws.Range("J132").Value = 7476.5
$ws.Range("K132").Value = 8370.900000000001
$ws.Range("L132").Value = 22429.5
$ws.Range("M132").Value = -5840.900000000001
$ws.Range("N132").Value = -27489.5
$ws.Range("H137").Value = 1841.4706
$ws.Range("I137").Value = 1742.8462
$ws.Range("K137").Value = 5228.5386
$ws.Range("M137").Value = -2678.5386

# ----- ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H25").Value = 1748
$ws.Range("I25").Value = 1748
$ws.Range("K25").Value = 1748
$ws.Range("M25").Value = -1346

# ----- BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2999.6667
$ws.Range("I20").Value = 2999.6667
$ws.Range("K20").Value = 2999.6667
$ws.Range("M20").Value = -2752.6667
$ws.Range("H33").Value = 9976
$ws.Range("I33").Value = 9976
$ws.Range("K33").Value = 9976
$ws.Range("M33").Value = -9640
$ws.Range("H99").Value = 4019.5833
$ws.Range("I99").Value = 4121.364
$ws.Range("K99").Value = 4121.364
$ws.Range("M99").Value = -2623.364
$ws.Range("H134").Value = 9627.875
$ws.Range("I134").Value = 5666
$ws.Range("K134").Value = 16998
$ws.Range("M134").Value = -14463

# ----- CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 500
$ws.Range("I7").Value = 500
$ws.Range("K7").Value = 500
$ws.Range("M7").Value = -387
$ws.Range("H39").Value = 0
$ws.Range("I39").Value = 0
$ws.Range("K39").Value = 0
$ws.Range("M39").ClearContents()
$ws.Range("H49").Value = 0
$ws.Range("I49").Value = 0
$ws.Range("K49").Value = 0
$ws.Range("M49").ClearContents()
$ws.Range("H59").Value = 39857.145
$ws.Range("J59").Value = 41800
$ws.Range("L59").Value = 41800
$ws.Range("N59").Value = -44090
$ws.Range("H132").Value = 2236
$ws.Range("I132").Value = 1322
$ws.Range("K132").Value = 3966
$ws.Range("M132").Value = -1436
$ws.Range("H134").Value = 2899.4
$ws.Range("I134").Value = 2499
$ws.Range("K134").Value = 7497
$ws.Range("M134").Value = -4962

# ----- CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1733025.1
$ws.Range("I4").Value = 1814472.5
$ws.Range("K4").Value = 5443417.5
$ws.Range("M4").Value = -5443305.5
$ws.Range("H11").Value = 500
$ws.Range("I11").Value = 500
$ws.Range("K11").Value = 1500
$ws.Range("M11").Value = -1360
$ws.Range("H12").Value = 30.333334
$ws.Range("J12").Value = 16.3
$ws.Range("L12").Value = 48.90000000000001
$ws.Range("N12").Value = -394.9
$ws.Range("H92").Value = 1353.8
$ws.Range("J92").Value = 390
$ws.Range("L92").Value = 1170
$ws.Range("N92").Value = -3666

# ----- GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 181
$ws.Range("I2").Value = 126.25
$ws.Range("K2").Value = 126.25
$ws.Range("M2").Value = -13.25
$ws.Range("H80").Value = 8600
$ws.Range("I80").Value = 8500
$ws.Range("J80").Value = 8666.666999999999
$ws.Range("K80").Value = 8500
$ws.Range("L80").Value = 8666.666999999999
$ws.Range("M80").Value = -7502
$ws.Range("N80").Value = -10662.667
$ws.Range("H83").Value = 8600
$ws.Range("I83").Value = 8500
$ws.Range("J83").Value = 8666.666999999999
$ws.Range("K83").Value = 42500
$ws.Range("L83").Value = 43333.335
$ws.Range("M83").Value = -37508
$ws.Range("N83").Value = -53317.335
$ws.Range("H113").Value = 875
$ws.Range("I113").Value = 875
$ws.Range("K113").Value = 875
$ws.Range("M113").Value = 1295

# ----- LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H32").Value = 19504.334
$ws.Range("I32").Value = 19504.334
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 19504.334
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -19187.334
$ws.Range("N32").ClearContents()
$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("M61").ClearContents()
$ws.Range("N61").ClearContents()
$ws.Range("H68").Value = 2734
$ws.Range("I68").Value = 2734
$ws.Range("K68").Value = 2734
$ws.Range("M68").Value = -1985
$ws.Range("H71").Value = 2734
$ws.Range("I71").Value = 2734
$ws.Range("K71").Value = 13670
$ws.Range("M71").Value = -9926
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("N113").ClearContents()
$ws.Range("H136").Value = 5469.857
$ws.Range("I136").Value = 8422.25
$ws.Range("K136").Value = 25266.75
$ws.Range("M136").Value = -22716.75

# ----- WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H33").Value = 12009
$ws.Range("I33").Value = 7019
$ws.Range("K33").Value = 7019
$ws.Range("M33").Value = -6769
$ws.Range("H34").Value = 11000
$ws.Range("I34").Value = 11000
$ws.Range("K34").Value = 11000
$ws.Range("M34").Value = -10797
$ws.Range("H36").Value = 12009
$ws.Range("I36").Value = 7019
$ws.Range("K36").Value = 7019
$ws.Range("M36").Value = -6769
$ws.Range("H81").Value = 1366.4
$ws.Range("J81").Value = 2885.5
$ws.Range("L81").Value = 5771
$ws.Range("N81").Value = -7893
$ws.Range("H84").Value = 1366.4
$ws.Range("J84").Value = 2885.5
$ws.Range("L84").Value = 28855
$ws.Range("N84").Value = -39463
$ws.Range("H107").Value = 439.6
$ws.Range("I107").Value = 266.33334
$ws.Range("K107").Value = 799.0000200000001
$ws.Range("M107").Value = 1120.99998
$ws.Range("H122").Value = 1960.85
$ws.Range("I122").Value = 1724.5294
$ws.Range("J122").Value = 3300
$ws.Range("K122").Value = 5173.5882
$ws.Range("L122").Value = 9900
$ws.Range("M122").Value = -2723.5882
$ws.Range("N122").Value = -14800
